# Update Work Week and Social Spending
# (Niger GDP per Capita series: refresh the "Data" sheet values for the
#  existing years and append the newly published years 2011-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$meta = $wb.Worksheets.Item("Metadata")

# New values for the already-present years (rows 2..62, years 1950..2010)
$newValues = @(
    983, 1004, 1023, 1042, 1063, 1082, 1103, 1122, 1143, 1162,
    1213, 1239, 1341, 1443, 1415, 1490, 1441, 1417, 1398, 1326,
    1337, 1382, 1282, 1039, 1105, 1050, 1033, 1087, 1205, 1261,
    1291, 1277, 1231, 1180, 956, 950, 969, 912, 939, 915,
    872, 865.913420964203, 785.263618323698, 771.541422741917,
    776.80252474935, 702.493394445968, 712.855845272195,
    693.10938249699, 755.247969163888, 736.917524763798,
    693.376547623471, 723.311358292809, 734.591580570774,
    758.114331368533, 725.131636267338, 757.243914388136,
    771.419659945294, 766.18694338525, 809.450392071507,
    775.253442693546, 810.157252880308
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $newValues[$i]
}

# Newly reported years, appended after the existing data block.
$newRows = @(
    @(2011, 799),
    @(2012, 864),
    @(2013, 880),
    @(2014, 911),
    @(2015, 913),
    @(2016, 925)
)

$startRow = 63
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 562
    $ws.Cells.Item($row, 2).Value = "Niger"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 5).Value = $newRows[$i][1]
}
